$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 165, shifting existing rows 165-205 down to 166-206.
$ws.Rows.Item(165).Insert()

# Populate the newly inserted row 165 with the new weekly record.
$ws.Range("A165").Value = 4
$ws.Range("B165").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C165").Value = "Los Lagos"
$ws.Range("D165").Value = 44508
$ws.Range("E165").Value = 10
$ws.Range("F165").Value = 100114014
$ws.Range("G165").Value = "Betarraga"
$ws.Range("H165").Value = "Sin especificar"
$ws.Range("I165").Value = "Primera"
$ws.Range("J165").Value = 500
$ws.Range("K165").Value = 1000
$ws.Range("L165").Value = 1000
$ws.Range("M165").Value = 1000
$ws.Range("N165").Value = "`$/paquete 5 unidades"
$ws.Range("O165").Value = "Región del Maule"
$ws.Range("P165").Value = 200
$ws.Range("Q165").Value = 5
$ws.Range("R165").Value = "Hortaliza"
